$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '61.415.36'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  -1.38%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.991.39'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  -0.32%  '
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '600.10'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  +3.25%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '143.77'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  -2.24%  '
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  +0.06%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.518'
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  -0.78%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '2.989.11'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  -0.42%  '
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  -1.04%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '6.05'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  +6.83%  '
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  +3.73%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '34.41'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  -0.78%  '
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  +2.34%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '3.488.77'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  -0.18%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '6.95'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  -0.89%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '61.382.09'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  -1.38%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '2.992.24'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  -0.29%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '450.34'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  -0.87%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '13.99'
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  +0.95%  '
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  +1.48%  '
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  +0.60%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '81.57'
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  +2.06%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '10.77'
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  +6.81%  '
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  -3.54%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '11.93'
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  -2.80%  '
$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  +0.24%  '
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  +3.15%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '7.18'
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  +0.29%  '
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  -1.64%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '27.24'
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  +1.22%  '
$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  +3.42%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.0₃0832'
$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  +5.07%  '
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  -0.28%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '5.77'
$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  +0.87%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '9.22'
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  +3.50%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '50.43'
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  +0.29%  '
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  -2.42%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.123'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  +10.62%  '
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  -1.05%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '397.69'
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  -3.13%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '39.90'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  +4.71%  '
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  +0.46%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.270'
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  -1.94%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.687.69'
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  -2.59%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '131.14'
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  +2.58%  '
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  -0.70%  '
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  +1.77%  '
